# ESN with 3 variables: insert a new "Time" column (C) before the existing
# "dc/dt" column, shifting dc/dt / prediction / (blank) / R-squared columns
# one place to the right (C->D, D->E, E->F, F->G).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column at C - this shifts C:F -> D:G and keeps each
#    shifted column's original per-cell styles/values intact.
$ws.Columns("C").Insert()

# 2) Header cell C1 = "Time", using the same "blank" style as the other
#    blank header cells (F1/G1), i.e. General format / general alignment.
$ws.Range("F1").Copy()
$ws.Range("C1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C1").Value = "Time"

# 3) Data cells C2:C8 = sequential Time index 1..7, with a dedicated
#    number format (#,##0) and general alignment - matching the new style
#    introduced for this column.
$ws.Range("F2").Copy()
$ws.Range("C2:C8").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C2:C8").NumberFormat = "#,##0"
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 3).Value = $r - 1
}

# 4) Remaining rows (9:145) in column C stay blank, matching style used for
#    the other blank columns.
$ws.Range("F10").Copy()
$ws.Range("C9:C145").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C9:C145").ClearContents()

# 5) Approximate the original column width/best-fit on the newly created
#    column (closest value reachable through the ColumnWidth property).
$ws.Columns(3).ColumnWidth = 12.6

$excel.CutCopyMode = $false
